$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voltage Table")

# --- Core data change: "updated voltage table for pac 12kv" ---
# The manual input cells (column H) of the 16/17/18-series block (rows 29-31)
# and the static reference block (rows 61-63) move from 16000 to 12000.
# Every other cell in those rows is a formula referencing column H (directly
# or indirectly), so they recalculate automatically.
$ws.Range("H29").Value = 12000
$ws.Range("H30").Value = 12000
$ws.Range("H31").Value = 12000

$ws.Range("H61").Value = 12000
$ws.Range("H62").Value = 12000
$ws.Range("H63").Value = 12000

# --- Formatting: row 31 picks up the same "section bottom" border/format
# already used on row 17 (thicker bottom border under cell H). Copy just the
# formats from H17 onto H31 so it reuses the existing style instead of
# creating a new one.
$ws.Range("H17").Copy() | Out-Null
$ws.Range("H31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- View state: active selection moved to H18 ---
$ws.Range("H18").Select() | Out-Null
